$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D column values remain text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '34.784.33'
$ws.Range('E2').Value = '  +2.74%  '
$ws.Range('D3').Value = '1.798.29'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '225.67'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').Value = '0.555'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '33.25'
$ws.Range('E8').Value = '  +7.77%  '
$ws.Range('D9').Value = '0.287'
$ws.Range('E9').Value = '  +3.03%  '
$ws.Range('D10').Value = '0.0674'
$ws.Range('E10').Value = '  +2.02%  '
$ws.Range('D11').Value = '0.0937'
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').Value = '2.064.21'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('D13').Value = '11.08'
$ws.Range('E13').Value = '  +10.83%  '
$ws.Range('D14').Value = '1.809.15'
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').Value = '0.638'
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('D16').Value = '34.782.78'
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('D17').Value = '4.31'
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('D18').Value = '69.54'
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('D19').Value = '256.83'
$ws.Range('E19').Value = '  +2.40%  '
$ws.Range('D20').Value = '0.0₃0766'
$ws.Range('E20').Value = '  +3.78%  '
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').Value = '10.48'
$ws.Range('E22').Value = '  +1.95%  '
$ws.Range('D23').Value = '4.25'
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('D24').Value = '2.14'
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('D25').Value = '159.04'
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('D26').Value = '16.49'
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').Value = '7.14'
$ws.Range('E27').Value = '  +3.26%  '
$ws.Range('D28').Value = '0.114'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '0.0521'
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '3.79'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').Value = '1.19'
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('D33').Value = '3.62'
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('D34').Value = '1.89'
$ws.Range('E34').Value = '  +8.50%  '
$ws.Range('D35').Value = '1.467.53'
$ws.Range('E35').Value = '  -1.37%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').Value = '1.06'
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.0191'
$ws.Range('E37').Value = '  +3.02%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '0.635'
$ws.Range('E38').Value = '  +1.50%  '
$ws.Range('D39').Value = '84.03'
$ws.Range('E39').Value = '  +1.44%  '
$ws.Range('D40').Value = '2.84'
$ws.Range('E40').Value = '  +5.64%  '
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('D42').Value = '0.904'
$ws.Range('E42').Value = '  +2.15%  '
$ws.Range('D43').Value = '2.10'
$ws.Range('E43').Value = '  +1.17%  '
$ws.Range('D44').Value = '0.0507'
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('D45').Value = '5.96'
$ws.Range('E45').Value = '  +4.19%  '
$ws.Range('D46').Value = '1.959.65'
$ws.Range('E46').Value = '  +1.23%  '
$ws.Range('E47').Value = '  -3.20%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '104.90'
$ws.Range('E48').Value = '  +6.56%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '12.03'
$ws.Range('E49').Value = '  +1.81%  '
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').Value = '49.79'
$ws.Range('E51').Value = '  -2.09%  '
